$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.32582968320327
$ws.Range("C2").Value = 8.472453906814785
$ws.Range("E2").Value = 12.47009151757337
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.672110882304977
$ws.Range("I2").Value = 26.0836635719414
$ws.Range("K2").Value = 8.747077556603333
$ws.Range("L2").Value = 10.19726921246547
$ws.Range("M2").Value = 13.76038602900983
$ws.Range("O2").Value = 25.96677269883756

$ws.Range("B3").Value = 11.07827039028585
$ws.Range("C3").Value = 8.451045178104968
$ws.Range("E3").Value = 12.50113403619955
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.673715595238764
$ws.Range("I3").Value = 26.20656583979776
$ws.Range("K3").Value = 8.566319071669923
$ws.Range("L3").Value = 10.20590862198843
$ws.Range("M3").Value = 13.72171817083736
$ws.Range("O3").Value = 26.08738321650708

$ws.Range("B4").Value = 10.92483694703489
$ws.Range("C4").Value = 8.43794495028714
$ws.Range("E4").Value = 12.52193474537765
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.674753052405734
$ws.Range("I4").Value = 26.28651275687345
$ws.Range("K4").Value = 8.454248197487992
$ws.Range("L4").Value = 10.21258808170198
$ws.Range("M4").Value = 13.69962973144842
$ws.Range("O4").Value = 26.1663325129804

$ws.Range("B5").Value = 10.86203685445207
$ws.Range("C5").Value = 8.432619279414224
$ws.Range("E5").Value = 12.53084916182748
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.675188979605722
$ws.Range("I5").Value = 26.32022069300565
$ws.Range("K5").Value = 8.408364342837283
$ws.Range("L5").Value = 10.21565627343974
$ws.Range("M5").Value = 13.69105033350608
$ws.Range("O5").Value = 26.19973576062754

$ws.Range("B6").Value = 10.85159490179191
$ws.Range("C6").Value = 8.431735771288764
$ws.Range("E6").Value = 12.53235585413043
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.675262160682138
$ws.Range("I6").Value = 26.32588609213031
$ws.Range("K6").Value = 8.400734185788348
$ws.Range("L6").Value = 10.21618667217774
$ws.Range("M6").Value = 13.68965137572163
$ws.Range("O6").Value = 26.20535668499758

$ws.Range("B7").Value = 10.92399100063063
$ws.Range("C7").Value = 8.437873072707255
$ws.Range("E7").Value = 12.5220531946402
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.674758878155882
$ws.Range("I7").Value = 26.28696278122834
$ws.Range("K7").Value = 8.453630179633594
$ws.Range("L7").Value = 10.21262805774682
$ws.Range("M7").Value = 13.69951231111952
$ws.Range("O7").Value = 26.16677801730501

$ws.Range("B8").Value = 11.24081833864349
$ws.Range("C8").Value = 8.465063782427194
$ws.Range("E8").Value = 12.48043392698687
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.672653386093049
$ws.Range("I8").Value = 26.12511029403635
$ws.Range("K8").Value = 8.685011652533872
$ws.Range("L8").Value = 10.19996308025958
$ws.Range("M8").Value = 13.74671378072484
$ws.Range("O8").Value = 26.00734354795464

$ws.Range("B9").Value = 11.84719176163476
$ws.Range("C9").Value = 8.518674764337437
$ws.Range("E9").Value = 12.41261494894678
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.668936549300904
$ws.Range("I9").Value = 25.84324180911181
$ws.Range("K9").Value = 9.127678698949451
$ws.Range("L9").Value = 10.18601006877534
$ws.Range("M9").Value = 13.85212319685135
$ws.Range("O9").Value = 25.73351767890398

$ws.Range("B10").Value = 12.27908893512715
$ws.Range("C10").Value = 8.558152815320391
$ws.Range("E10").Value = 12.37117999473905
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.666454425632965
$ws.Range("I10").Value = 25.65772782130815
$ws.Range("K10").Value = 9.443028558427672
$ws.Range("L10").Value = 10.18235633612434
$ws.Range("M10").Value = 13.93701411277351
$ws.Range("O10").Value = 25.55599704598211

$ws.Range("B11").Value = 12.47174788322016
$ws.Range("C11").Value = 8.576113202426679
$ws.Range("E11").Value = 12.35414806738059
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.665378690767454
$ws.Range("I11").Value = 25.57800096345482
$ws.Range("K11").Value = 9.583738020599588
$ws.Range("L11").Value = 10.18211751859478
$ws.Range("M11").Value = 13.97716193205524
$ws.Range("O11").Value = 25.48037479278669

$ws.Range("B12").Value = 12.54408833902397
$ws.Range("C12").Value = 8.582912900556282
$ws.Range("E12").Value = 12.34795945065111
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.664978975351115
$ws.Range("I12").Value = 25.54848014917539
$ws.Range("K12").Value = 9.636579823096056
$ws.Range("L12").Value = 10.18223083081101
$ws.Range("M12").Value = 13.99257701569509
$ws.Range("O12").Value = 25.45247689143034

$ws.Range("B13").Value = 12.52853690904364
$ws.Range("C13").Value = 8.581448559026844
$ws.Range("E13").Value = 12.34928067645116
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.665064721909563
$ws.Range("I13").Value = 25.55480819639061
$ws.Range("K13").Value = 9.62521976217745
$ws.Range("L13").Value = 10.18219738073026
$ws.Range("M13").Value = 13.98924780577137
$ws.Range("O13").Value = 25.45845234515236

$ws.Range("B14").Value = 12.47771207435136
$ws.Range("C14").Value = 8.576672658005471
$ws.Range("E14").Value = 12.35363369692446
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.665345652988316
$ws.Range("I14").Value = 25.57555884407899
$ws.Range("K14").Value = 9.588094460872437
$ws.Range("L14").Value = 10.18212276299544
$ws.Range("M14").Value = 13.97842593930959
$ws.Range("O14").Value = 25.47806480974584

$ws.Range("B15").Value = 12.44649835835356
$ws.Range("C15").Value = 8.573747031964189
$ws.Range("E15").Value = 12.35633402943362
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.665518725534486
$ws.Range("I15").Value = 25.58835645054299
$ws.Range("K15").Value = 9.565295259483179
$ws.Range("L15").Value = 10.18210356280746
$ws.Range("M15").Value = 13.97182459408308
$ws.Range("O15").Value = 25.49017421508263

$ws.Range("B16").Value = 12.26641565433743
$ws.Range("C16").Value = 8.556978945725223
$ws.Range("E16").Value = 12.37232961089284
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.666525798791692
$ws.Range("I16").Value = 25.66303195584232
$ws.Range("K16").Value = 9.433773532397243
$ws.Range("L16").Value = 10.18240050667081
$ws.Range("M16").Value = 13.93442046883005
$ws.Range("O16").Value = 25.56104246560523

$ws.Range("B17").Value = 12.15491492801122
$ws.Range("C17").Value = 8.546691447955354
$ws.Range("E17").Value = 12.38260757332342
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.667157255667856
$ws.Range("I17").Value = 25.71003709364063
$ws.Range("K17").Value = 9.352351686907099
$ws.Range("L17").Value = 10.18294657922191
$ws.Range("M17").Value = 13.91186056329985
$ws.Range("O17").Value = 25.60583281909591

$ws.Range("B18").Value = 12.09042744223729
$ws.Range("C18").Value = 8.540774609567753
$ws.Range("E18").Value = 12.38869023202895
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.667525480956547
$ws.Range("I18").Value = 25.73751227912015
$ws.Range("K18").Value = 9.305264263867675
$ws.Range("L18").Value = 10.18339469169805
$ws.Range("M18").Value = 13.89902926400446
$ws.Range("O18").Value = 25.63207813210918

$ws.Range("B19").Value = 12.0685341972089
$ws.Range("C19").Value = 8.538771366984086
$ws.Range("E19").Value = 12.39077910152798
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.667651020382825
$ws.Range("I19").Value = 25.74689033886215
$ws.Range("K19").Value = 9.289278836756285
$ws.Range("L19").Value = 10.18356946074285
$ws.Range("M19").Value = 13.89470987802519
$ws.Range("O19").Value = 25.64104729983307

$ws.Range("B20").Value = 12.16682163451303
$ws.Range("C20").Value = 8.547786555516838
$ws.Range("E20").Value = 12.38149576600619
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.667089515889801
$ws.Range("I20").Value = 25.70498787937741
$ws.Range("K20").Value = 9.361045999883167
$ws.Range("L20").Value = 10.18287458410057
$ws.Range("M20").Value = 13.91424720010053
$ws.Range("O20").Value = 25.60101480829466

$ws.Range("B21").Value = 12.49265776961381
$ws.Range("C21").Value = 8.578075510881913
$ws.Range("E21").Value = 12.35234802805291
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.665262929604832
$ws.Range("I21").Value = 25.56944569723361
$ws.Range("K21").Value = 9.599011412496909
$ws.Range("L21").Value = 10.18213915809381
$ws.Range("M21").Value = 13.98159889674596
$ws.Range("O21").Value = 25.47228410355251

$ws.Range("B22").Value = 12.70199634096565
$ws.Range("C22").Value = 8.597861469063503
$ws.Range("E22").Value = 12.33481948642346
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.664113676897816
$ws.Range("I22").Value = 25.48476612110163
$ws.Range("K22").Value = 9.751940687253082
$ws.Range("L22").Value = 10.18284562942494
$ws.Range("M22").Value = 14.02684874453646
$ws.Range("O22").Value = 25.39245658026831

$ws.Range("B23").Value = 12.59061999571509
$ws.Range("C23").Value = 8.587302791651309
$ws.Range("E23").Value = 12.34403570821782
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.664722992351833
$ws.Range("I23").Value = 25.52960411243003
$ws.Range("K23").Value = 9.670571639391364
$ws.Range("L23").Value = 10.18236027459602
$ws.Range("M23").Value = 14.00258809970362
$ws.Range("O23").Value = 25.43466789850331

$ws.Range("B24").Value = 12.16143980416406
$ws.Range("C24").Value = 8.547291464798839
$ws.Range("E24").Value = 12.38199787295612
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.667120124866941
$ws.Range("I24").Value = 25.70726922329569
$ws.Range("K24").Value = 9.357116159913675
$ws.Range("L24").Value = 10.1829067151099
$ws.Range("M24").Value = 13.91316776898325
$ws.Range("O24").Value = 25.60319148991763

$ws.Range("B25").Value = 11.68523113120315
$ws.Range("C25").Value = 8.504149393825296
$ws.Range("E25").Value = 12.4294866935908
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.669898206170195
$ws.Range("I25").Value = 25.9157002532721
$ws.Range("K25").Value = 9.009444352462415
$ws.Range("L25").Value = 10.18862317035059
$ws.Range("M25").Value = 13.82227031603931
$ws.Range("O25").Value = 25.80343955816626
